# Auto-generated edit script applying numeric corrections to the
# per-profession Leve profit tables (H/I/J/K/L/M/N columns), as
# produced by the scheduled price-refresh runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 1195.55
$ws.Range("I92").Value = 366.06668
$ws.Range("K92").Value = 366.06668
$ws.Range("M92").Value = 881.93332
$ws.Range("H98").Value = 1297.3334
$ws.Range("J98").Value = 2150
$ws.Range("L98").Value = 2150
$ws.Range("N98").Value = -5146
$ws.Range("H122").Value = 1297.3334
$ws.Range("J122").Value = 2150
$ws.Range("L122").Value = 6450
$ws.Range("N122").Value = -11350
$ws.Range("H135").Value = 1857.3889
$ws.Range("I135").Value = 790.17645
$ws.Range("K135").Value = 7111.58805
$ws.Range("M135").Value = -4576.58805
$ws.Range("H138").Value = 3437.5217
$ws.Range("I138").Value = 4366.5
$ws.Range("K138").Value = 13099.5
$ws.Range("M138").Value = -7959.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1112.3462
$ws.Range("I32").Value = 996.7755
$ws.Range("J32").Value = 3000
$ws.Range("K32").Value = 996.7755
$ws.Range("L32").Value = 3000
$ws.Range("M32").Value = -709.7755
$ws.Range("N32").Value = -3574
$ws.Range("H45").Value = 2955.2104
$ws.Range("I45").Value = 3044.353
$ws.Range("K45").Value = 3044.353
$ws.Range("M45").Value = -2667.353
$ws.Range("H102").Value = 33872.273
$ws.Range("I102").Value = 36959.5
$ws.Range("J102").Value = 3000
$ws.Range("K102").Value = 36959.5
$ws.Range("L102").Value = 3000
$ws.Range("M102").Value = -35337.5
$ws.Range("N102").Value = -6244
$ws.Range("H132").Value = 2286.48
$ws.Range("I132").Value = 1567.5476
$ws.Range("K132").Value = 4702.642800000001
$ws.Range("M132").Value = -2172.642800000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1040.95
$ws.Range("I20").Value = 852.75
$ws.Range("J20").Value = 1323.25
$ws.Range("K20").Value = 852.75
$ws.Range("L20").Value = 1323.25
$ws.Range("M20").Value = -605.75
$ws.Range("N20").Value = -1817.25
$ws.Range("H94").Value = 720.2353000000001
$ws.Range("I94").Value = 518.4286
$ws.Range("J94").Value = 1046.2307
$ws.Range("K94").Value = 518.4286
$ws.Range("L94").Value = 1046.2307
$ws.Range("M94").Value = -67.42859999999996
$ws.Range("N94").Value = -1948.2307
$ws.Range("H99").Value = 9231.294
$ws.Range("I99").Value = 3787.3572
$ws.Range("K99").Value = 3787.3572
$ws.Range("M99").Value = -2289.3572
$ws.Range("H134").Value = 3726.8655
$ws.Range("I134").Value = 1698.0555
$ws.Range("K134").Value = 5094.166499999999
$ws.Range("M134").Value = -2559.166499999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value = 932
$ws.Range("I25").Value = 932
$ws.Range("K25").Value = 932
$ws.Range("M25").Value = -758
$ws.Range("H31").Value = 2407
$ws.Range("I31").Value = 1141.8108
$ws.Range("J31").Value = 4747.6
$ws.Range("K31").Value = 1141.8108
$ws.Range("L31").Value = 4747.6
$ws.Range("M31").Value = -846.8108
$ws.Range("N31").Value = -5337.6
$ws.Range("H34").Value = 2407
$ws.Range("I34").Value = 1141.8108
$ws.Range("J34").Value = 4747.6
$ws.Range("K34").Value = 1141.8108
$ws.Range("L34").Value = 4747.6
$ws.Range("M34").Value = -939.8108
$ws.Range("N34").Value = -5151.6
$ws.Range("H58").Value = 2815.7778
$ws.Range("J58").Value = 4050.3572
$ws.Range("L58").Value = 4050.3572
$ws.Range("N58").Value = -4456.3572
$ws.Range("H132").Value = 40003910
$ws.Range("I132").Value = 60609252
$ws.Range("J132").Value = 5311
$ws.Range("K132").Value = 181827756
$ws.Range("L132").Value = 15933
$ws.Range("M132").Value = -181825226
$ws.Range("N132").Value = -20993
$ws.Range("H134").Value = 2678.6316
$ws.Range("I134").Value = 2177.1428
$ws.Range("K134").Value = 6531.428400000001
$ws.Range("M134").Value = -3996.428400000001
$ws.Range("H136").Value = 2815.7778
$ws.Range("J136").Value = 4050.3572
$ws.Range("L136").Value = 12151.0716
$ws.Range("N136").Value = -17251.0716

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2289.6
$ws.Range("I5").Value = 482.66666
$ws.Range("J5").Value = 5000
$ws.Range("K5").Value = 1447.99998
$ws.Range("L5").Value = 15000
$ws.Range("M5").Value = -1335.99998
$ws.Range("N5").Value = -15224
$ws.Range("H68").Value = 11501.454
$ws.Range("J68").Value = 12401.6
$ws.Range("L68").Value = 37204.8
$ws.Range("N68").Value = -38826.8
$ws.Range("H71").Value = 11501.454
$ws.Range("J71").Value = 12401.6
$ws.Range("L71").Value = 111614.4
$ws.Range("N71").Value = -119726.4
$ws.Range("H135").Value = 2289.6
$ws.Range("I135").Value = 482.66666
$ws.Range("J135").Value = 5000
$ws.Range("K135").Value = 4343.99994
$ws.Range("L135").Value = 45000
$ws.Range("M135").Value = -1808.99994
$ws.Range("N135").Value = -50070

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H74").Value = 91754
$ws.Range("I74").Value = 75000
$ws.Range("K74").Value = 75000
$ws.Range("M74").Value = -74064
$ws.Range("H77").Value = 91754
$ws.Range("I77").Value = 75000
$ws.Range("K77").Value = 225000
$ws.Range("M77").Value = -220320
$ws.Range("H113").Value = 1912.6923
$ws.Range("I113").Value = 1671.5555
$ws.Range("J113").Value = 2455.25
$ws.Range("K113").Value = 1671.5555
$ws.Range("L113").Value = 2455.25
$ws.Range("M113").Value = 498.4445000000001
$ws.Range("N113").Value = -6795.25
$ws.Range("H132").Value = 12993243
$ws.Range("I132").Value = 16399030
$ws.Range("J132").Value = 8682.625
$ws.Range("K132").Value = 49197090
$ws.Range("L132").Value = 26047.875
$ws.Range("M132").Value = -49194560
$ws.Range("N132").Value = -31107.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1198.2222
$ws.Range("I22").Value = 398.75
$ws.Range("J22").Value = 1837.8
$ws.Range("K22").Value = 398.75
$ws.Range("L22").Value = 1837.8
$ws.Range("M22").Value = -103.75
$ws.Range("N22").Value = -2427.8
$ws.Range("H27").Value = 1198.2222
$ws.Range("I27").Value = 398.75
$ws.Range("J27").Value = 1837.8
$ws.Range("K27").Value = 398.75
$ws.Range("L27").Value = 1837.8
$ws.Range("M27").Value = -291.75
$ws.Range("N27").Value = -2051.8
$ws.Range("H124").Value = 54166.668
$ws.Range("J124").Value = 54166.668
$ws.Range("L124").Value = 54166.668
$ws.Range("N124").Value = -63986.668
$ws.Range("H132").Value = 1664.4445
$ws.Range("I132").Value = 1739.3914
$ws.Range("J132").Value = 1233.5
$ws.Range("K132").Value = 5218.174199999999
$ws.Range("L132").Value = 3700.5
$ws.Range("M132").Value = -2688.174199999999
$ws.Range("N132").Value = -8760.5
$ws.Range("H136").Value = 5122.15
$ws.Range("I136").Value = 1911.25
$ws.Range("J136").Value = 7262.75
$ws.Range("K136").Value = 5733.75
$ws.Range("L136").Value = 21788.25
$ws.Range("M136").Value = -3183.75
$ws.Range("N136").Value = -26888.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1636.4375
$ws.Range("I100").Value = 1158
$ws.Range("K100").Value = 2316
$ws.Range("M100").Value = -1775
$ws.Range("H132").Value = 7696078
$ws.Range("I132").Value = 10528475
$ws.Range("K132").Value = 31585425
$ws.Range("M132").Value = -31582895
$ws.Range("H133").Value = 89999.664
$ws.Range("J133").Value = 89999.664
$ws.Range("L133").Value = 89999.664
$ws.Range("N133").Value = -100119.664
$ws.Range("H136").Value = 23816588
$ws.Range("I136").Value = 41671904
$ws.Range("K136").Value = 125015712
$ws.Range("M136").Value = -125013162
